# Re-sort the sheet tabs: "总计" (summary) should come first, followed by
# "2020-Q4" (the detail/fund table). Content of each sheet is left untouched;
# only the tab order changes.
$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ4    = $wb.Worksheets.Item("2020-Q4")

# Move "总计" so it sits right before "2020-Q4" -> final order: 总计, 2020-Q4
$wsTotal.Move($wsQ4)
